# Insert a new data row at row 8 (pushes the existing rows 8..85 down to 9..86),
# and populate it with the new "Poroto granado" price record
# (Fecha 2022-12-26, Volumen 155, Precios 50000, Origen "Región del Maule", Precio $/Kg 2000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value()  = 10
$ws.Cells.Item(8, 2).Value()  = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value()  = "La Araucanía"
$ws.Cells.Item(8, 4).Value()  = 44921
$ws.Cells.Item(8, 5).Value()  = 9
$ws.Cells.Item(8, 6).Value()  = 100112030
$ws.Cells.Item(8, 7).Value()  = "Poroto granado"
$ws.Cells.Item(8, 8).Value()  = "Sin especificar"
$ws.Cells.Item(8, 9).Value()  = "Primera"
$ws.Cells.Item(8, 10).Value() = 155
$ws.Cells.Item(8, 11).Value() = 50000
$ws.Cells.Item(8, 12).Value() = 50000
$ws.Cells.Item(8, 13).Value() = 50000
$ws.Cells.Item(8, 14).Value() = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value() = "Región del Maule"
$ws.Cells.Item(8, 16).Value() = 2000
$ws.Cells.Item(8, 17).Value() = 25
$ws.Cells.Item(8, 18).Value() = "Hortaliza"
